# Modified sharp sensor price
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sharp IR sensor (GP2Y0A21YK0F) unit price changed
$ws.Range("D2").Value = 65.08

# Nema17 stepper motor (CSGSHJ) unit price changed
$ws.Range("D3").Value = 45.77

# A new (otherwise blank) cell with a single space value appears further down the sheet
$ws.Range("F22").Value = " "

# Update the active selection shown in the sheet view
[void]$ws.Range("N8").Select()
